$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    4944.292321154038,
    4887.559847582665,
    4825.545484664752,
    4794.057447783162,
    4794.010622970725,
    4890.506649537315,
    5267.449152182856,
    6573.678702600464,
    8191.789500000001,
    14018.138,
    14986.9245,
    14980.626,
    15169.25375,
    15666.69225,
    15813.73925,
    16092.652,
    16268.66525,
    16049.11175,
    14955.2725,
    13776.544,
    12164.7975,
    9625.974749999999,
    6491.415499999999,
    5924.492750000001
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
